# Actualización automática 2025-10-21 15:30:08
#
# Updates three recorded sales figures for GUERRERO FAREZ FABIAN MAURICIO
# (advisor) and propagates the resulting deltas through the monthly
# totals / compliance summary sheets that cache these aggregates as
# static values.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "VENTAS POR GRUPO" (sales broken down by product group)
# ---------------------------------------------------------------
$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")

# FEIJOO MARIN MAURICIO ENRIQUE - PORCELANATO
$wsGrupo.Range("M24").Value = 5320.31

# ORTEGA ROMAN KLEBER ERWIN - PANELES DECORATIVOS
$wsGrupo.Range("K36").Value = 1446.48

# SALDARRIAGA ECHEVERRIA BRYAN STEVEN - FREGADEROS DE COCINA
$wsGrupo.Range("E50").Value = 55.65

# Running "N de 54" tally for the FREGADEROS DE COCINA column now that
# E50 became non-zero.
$wsGrupo.Range("E56").Value = "4 de 54"

# ---------------------------------------------------------------
# Sheet "VENTA MENSUAL" (monthly sales - "octubre" column)
# ---------------------------------------------------------------
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")

$wsMensual.Range("F24").Value = 6106.91
$wsMensual.Range("F36").Value = 10980.93
$wsMensual.Range("F50").Value = 55.65

# Column total for "octubre"
$wsMensual.Range("F60").Value = 52255.17

# ---------------------------------------------------------------
# Sheet "CUMPLIMIENTO MENSUAL" (monthly compliance by product group)
# ---------------------------------------------------------------
$wsCumpl = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# Column D widened slightly to fit the larger VENTA figures.
$wsCumpl.Columns.Item(4).ColumnWidth = 13.1

# FREGADEROS DE COCINA row
$wsCumpl.Range("D4").Value = 1375.03
$wsCumpl.Range("E4").Value = -331.80711473472
$wsCumpl.Range("F4").Value = 1.318059658603391

# PANELES DECORATIVOS row
$wsCumpl.Range("D10").Value = 5268.17
$wsCumpl.Range("E10").Value = -1387.09016465608
$wsCumpl.Range("F10").Value = 1.357397998367422

# PORCELANATO row
$wsCumpl.Range("D12").Value = 26866.39
$wsCumpl.Range("E12").Value = 25796.73
$wsCumpl.Range("F12").Value = 0.5101556839017514

# TOTAL row
$wsCumpl.Range("D14").Value = 50160.68000000001
$wsCumpl.Range("E14").Value = 48855.82661190613
$wsCumpl.Range("F14").Value = 0.5065890700083383
